$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Orig Pwd" value for CarlG@TestIncidentQueue.onmicrosoft.com (row 5)
# was changed from the old generated password "Pufa7292" to "P@ssw0rd4".
$ws.Range("C5").Value = "P@ssw0rd4"

# That same cell picked up a hyperlink as part of the Worker Queue Screen
# work (e.g. linking to the reset-password reference). Add it, then put
# the cell's formatting back the way it was (Excel's default Hyperlinks.Add
# call re-styles the cell with the built-in "Hyperlink" style, but the
# target cell here keeps its plain/default look).
$ws.Hyperlinks.Add($ws.Range("C5"), "https://support.microsoft.com/")
$ws.Range("C5").Style = "Normal"

# Selection left on D5 when the file was saved.
$ws.Range("D5").Select()
